$p = $ppt.ActivePresentation

# Slide 10 (title "Cnclusione" -> "Cnclusioni"): simple typo fix, one extra letter "i"
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Cnclusioni"

# Slide 2 (title "Idea di base" -> "Obbiettivi"): replace the whole title text, then
# split off the leading "O" into its own run tagged as English (en-US) while the
# remainder ("bbiettivi") keeps the original Italian (it-IT) run formatting.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Text = "bbiettivi"
$tr2.InsertBefore("O")
$firstLetter = $tr2.Characters(1, 1)
$firstLetter.LanguageID = "en-US"
